{"js": "// The document contains the markup `<id>p123r_1</id>` split across three\n// separate runs (the `<id>` and `</id>` tags carry Courier-New / brown\n// formatting, while `p123r_1` sits in an unformatted run in between).\n// The edit collapses those three runs into a single run whose text is\n// the full `<id>p123r_1</id>` string, keeping the first run's formatting.\n//\n// A Word `search()` match can span multiple runs, and replacing its text\n// merges the matched runs into one (taking the formatting of the first\n// run in the match) - exactly the behaviour we need here.\nconst results = context.document.body.search(\"<id>p123r_1</id>\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find '<id>p123r_1</id>' in the document body\");\n}\n\nconst target = results.items[0];\ntarget.insertText(\"<id>p123r_1</id>\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The document contains the markup `<id>p123r_1</id>` split across three\n# separate runs (the `<id>` and `</id>` tags carry Courier-New / brown\n# formatting, while `p123r_1` sits in an unformatted run in between).\n# The edit collapses those three runs into a single run whose text is\n# the full `<id>p123r_1</id>` string, keeping the first run's formatting.\n#\n# Word's Find/Replace can match text spanning multiple runs; replacing it\n# merges the matched runs into one (taking the formatting of the first\n# run in the match) - exactly the behaviour we need here.\n$d = $word.ActiveDocument\n$rng = $d.Content\n$found = $rng.Find.Execute(\"<id>p123r_1</id>\", $false, $false, $false, $false, $false, $true, 1, $false, \"<id>p123r_1</id>\", 2)\n"}
